# Shop Page Dev Completed
# Append a new "Product" / "Blackberry" row to the data set on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Product"
$ws.Range("B7").Value = "Blackberry"
